$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-date the RC 300 / RC 350 rows (rows 2-9) from MY2020 to MY2021 and
#    bump their base MSRP values.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 2021
$ws.Range("D2").Value = 42120

$ws.Range("C3").Value = 2021
$ws.Range("D3").Value = 46590

$ws.Range("C4").Value = 2021
$ws.Range("D4").Value = 44810

$ws.Range("C5").Value = 2021
$ws.Range("D5").Value = 48765

$ws.Range("C6").Value = 2021
$ws.Range("D6").Value = 45050

$ws.Range("C7").Value = 2021
$ws.Range("D7").Value = 49520

$ws.Range("C8").Value = 2021
$ws.Range("D8").Value = 47215

$ws.Range("C9").Value = 2021
$ws.Range("D9").Value = 51130

# ---------------------------------------------------------------------------
# 2. Re-date / re-price the LS 500 Inspiration Series rows (53-54).
# ---------------------------------------------------------------------------
$ws.Range("C53").Value = 2021
$ws.Range("D53").Value = 65875

$ws.Range("C54").Value = 2021
$ws.Range("D54").Value = 96675

# ---------------------------------------------------------------------------
# 3. Append the four new RC "Black Line" trims (rows 95-98). New shared
#    strings must be introduced in trim-code order first, then name order,
#    matching the authoring tool's write sequence.
# ---------------------------------------------------------------------------
$ws.Range("A95").Value = "9203SE"
$ws.Range("A96").Value = "9207SE"
$ws.Range("A97").Value = "9213SE"
$ws.Range("A98").Value = "9217SE"

$ws.Range("B95").Value = "RC 300 F SPORT Black Line"
$ws.Range("B96").Value = "RC 300 AWD F SPORT Black Line"
$ws.Range("B97").Value = "RC 350 F SPORT Black Line"
$ws.Range("B98").Value = "RC 350 AWD F SPORT Black Line"

$ws.Range("C95").Value = 2021
$ws.Range("C96").Value = 2021
$ws.Range("C97").Value = 2021
$ws.Range("C98").Value = 2021

$ws.Range("D95").Value = 48735
$ws.Range("D96").Value = 50910
$ws.Range("D97").Value = 51665
$ws.Range("D98").Value = 53275
$ws.Range("D95:D98").NumberFormat = "$#,##0_);[Red]($#,##0)"

$ws.Range("E95").Value = 1025
$ws.Range("E96").Value = 1025
$ws.Range("E97").Value = 1025
$ws.Range("E98").Value = 1025
$ws.Range("E95:E98").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# ---------------------------------------------------------------------------
# 4. Update the view state: scroll the window so row 40 is at the top and
#    move the active selection to C55 (mirrors the author's on-screen state
#    when the workbook was saved).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("C55").Select()
